$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1928310.4  # ALC!H17 (was 2005365.4)
$ws.Cells.Item(17, 10).Value = 1928310.4  # ALC!J17 (was 2005365.4)
$ws.Cells.Item(17, 12).Value = 5784931.199999999  # ALC!L17 (was 6016096.199999999)
$ws.Cells.Item(17, 14).Value = -5785267.199999999  # ALC!N17 (was -6016432.199999999)

$ws.Cells.Item(40, 8).Value = 2040.4  # ALC!H40 (was 2610.7693)
$ws.Cells.Item(40, 9).Value = 2100  # ALC!I40 (was 2605.4546)
$ws.Cells.Item(40, 10).Value = 1951  # ALC!J40 (was 2640)
$ws.Cells.Item(40, 11).Value = 2100  # ALC!K40 (was 2605.4546)
$ws.Cells.Item(40, 12).Value = 1951  # ALC!L40 (was 2640)
$ws.Cells.Item(40, 13).Value = -1925  # ALC!M40 (was -2430.4546)
$ws.Cells.Item(40, 14).Value = -2301  # ALC!N40 (was -2990)

$ws.Cells.Item(58, 8).Value = 404.16666  # ALC!H58 (was 257.5)
$ws.Cells.Item(58, 9).Value = 106.25  # ALC!I58 (was 109)
$ws.Cells.Item(58, 11).Value = 318.75  # ALC!K58 (was 327)
$ws.Cells.Item(58, 13).Value = -168.75  # ALC!M58 (was -177)

$ws.Cells.Item(62, 8).Value = 3561.8  # ALC!H62 (was 3724.5)
$ws.Cells.Item(62, 9).Value = 2530  # ALC!I62 (was 2515.3845)
$ws.Cells.Item(62, 10).Value = 4593.6  # ALC!J62 (was 5970)
$ws.Cells.Item(62, 11).Value = 2530  # ALC!K62 (was 2515.3845)
$ws.Cells.Item(62, 12).Value = 4593.6  # ALC!L62 (was 5970)
$ws.Cells.Item(62, 13).Value = -1906  # ALC!M62 (was -1891.3845)
$ws.Cells.Item(62, 14).Value = -5841.6  # ALC!N62 (was -7218)

$ws.Cells.Item(64, 8).Value = 4747.7  # ALC!H64 (was 3949.1091)
$ws.Cells.Item(64, 9).Value = 3565  # ALC!I64 (was 3279.5454)
$ws.Cells.Item(64, 10).Value = 8000.125  # ALC!J64 (was 6627.364)
$ws.Cells.Item(64, 11).Value = 3565  # ALC!K64 (was 3279.5454)
$ws.Cells.Item(64, 12).Value = 8000.125  # ALC!L64 (was 6627.364)
$ws.Cells.Item(64, 13).Value = -3317  # ALC!M64 (was -3031.5454)
$ws.Cells.Item(64, 14).Value = -8496.125  # ALC!N64 (was -7123.364)

$ws.Cells.Item(65, 8).Value = 3561.8  # ALC!H65 (was 3724.5)
$ws.Cells.Item(65, 9).Value = 2530  # ALC!I65 (was 2515.3845)
$ws.Cells.Item(65, 10).Value = 4593.6  # ALC!J65 (was 5970)
$ws.Cells.Item(65, 11).Value = 12650  # ALC!K65 (was 12576.9225)
$ws.Cells.Item(65, 12).Value = 22968  # ALC!L65 (was 29850)
$ws.Cells.Item(65, 13).Value = -9530  # ALC!M65 (was -9456.922500000001)
$ws.Cells.Item(65, 14).Value = -29208  # ALC!N65 (was -36090)

$ws.Cells.Item(67, 8).Value = 4747.7  # ALC!H67 (was 3949.1091)
$ws.Cells.Item(67, 9).Value = 3565  # ALC!I67 (was 3279.5454)
$ws.Cells.Item(67, 10).Value = 8000.125  # ALC!J67 (was 6627.364)
$ws.Cells.Item(67, 11).Value = 3565  # ALC!K67 (was 3279.5454)
$ws.Cells.Item(67, 12).Value = 8000.125  # ALC!L67 (was 6627.364)
$ws.Cells.Item(67, 13).Value = -2707  # ALC!M67 (was -2421.5454)
$ws.Cells.Item(67, 14).Value = -9716.125  # ALC!N67 (was -8343.364)

$ws.Cells.Item(74, 8).Value = 4500  # ALC!H74 (was 4923.077)
$ws.Cells.Item(74, 9).Value = 0  # ALC!I74 (was 4666.6665)
$ws.Cells.Item(74, 10).Value = 4500  # ALC!J74 (was 5000)
$ws.Cells.Item(74, 11).Value = 0  # ALC!K74 (was 4666.6665)
$ws.Cells.Item(74, 12).Value = 4500  # ALC!L74 (was 5000)
$ws.Cells.Item(74, 13).Value = $null  # ALC!M74 (was -3730.6665)
$ws.Cells.Item(74, 14).Value = -6372  # ALC!N74 (was -6872)

$ws.Cells.Item(77, 8).Value = 4500  # ALC!H77 (was 4923.077)
$ws.Cells.Item(77, 9).Value = 0  # ALC!I77 (was 4666.6665)
$ws.Cells.Item(77, 10).Value = 4500  # ALC!J77 (was 5000)
$ws.Cells.Item(77, 11).Value = 0  # ALC!K77 (was 23333.3325)
$ws.Cells.Item(77, 12).Value = 22500  # ALC!L77 (was 25000)
$ws.Cells.Item(77, 13).Value = $null  # ALC!M77 (was -18653.3325)
$ws.Cells.Item(77, 14).Value = -31860  # ALC!N77 (was -34360)

$ws.Cells.Item(112, 8).Value = 1229.4783  # ALC!H112 (was 1218.5652)
$ws.Cells.Item(112, 10).Value = 1283.9  # ALC!J112 (was 1271.35)
$ws.Cells.Item(112, 12).Value = 3851.7  # ALC!L112 (was 3814.05)
$ws.Cells.Item(112, 14).Value = -6067.700000000001  # ALC!N112 (was -6030.049999999999)

$ws.Cells.Item(113, 8).Value = 8071.4287  # ALC!H113 (was 7889.6553)
$ws.Cells.Item(113, 10).Value = 10242.263  # ALC!J113 (was 9870.15)
$ws.Cells.Item(113, 12).Value = 10242.263  # ALC!L113 (was 9870.15)
$ws.Cells.Item(113, 14).Value = -16750.263  # ALC!N113 (was -16378.15)

$ws.Cells.Item(115, 8).Value = 852.5  # ALC!H115 (was 1145.8334)
$ws.Cells.Item(115, 9).Value = 305.83334  # ALC!I115 (was 375)
$ws.Cells.Item(115, 10).Value = 1262.5  # ALC!J115 (was 1300)
$ws.Cells.Item(115, 11).Value = 917.5000200000001  # ALC!K115 (was 1125)
$ws.Cells.Item(115, 12).Value = 3787.5  # ALC!L115 (was 3900)
$ws.Cells.Item(115, 13).Value = 649.4999799999999  # ALC!M115 (was 442)
$ws.Cells.Item(115, 14).Value = -6921.5  # ALC!N115 (was -7034)

$ws.Cells.Item(118, 8).Value = 860.37036  # ALC!H118 (was 904.5833)
$ws.Cells.Item(118, 9).Value = 320.83334  # ALC!I118 (was 353)
$ws.Cells.Item(118, 10).Value = 1292  # ALC!J118 (was 1298.5714)
$ws.Cells.Item(118, 11).Value = 962.5000200000001  # ALC!K118 (was 1059)
$ws.Cells.Item(118, 12).Value = 3876  # ALC!L118 (was 3895.7142)
$ws.Cells.Item(118, 13).Value = 694.4999799999999  # ALC!M118 (was 598)
$ws.Cells.Item(118, 14).Value = -7190  # ALC!N118 (was -7209.7142)

$ws.Cells.Item(129, 8).Value = 916.34424  # ALC!H129 (was 964.75)
$ws.Cells.Item(129, 9).Value = 376.07144  # ALC!I129 (was 423)
$ws.Cells.Item(129, 10).Value = 1077.2766  # ALC!J129 (was 1167.9062)
$ws.Cells.Item(129, 11).Value = 1128.21432  # ALC!K129 (was 1269)
$ws.Cells.Item(129, 12).Value = 3231.8298  # ALC!L129 (was 3503.7186)
$ws.Cells.Item(129, 13).Value = 3871.78568  # ALC!M129 (was 3731)
$ws.Cells.Item(129, 14).Value = -13231.8298  # ALC!N129 (was -13503.7186)

$ws.Cells.Item(138, 8).Value = 3383.93  # ALC!H138 (was 2494.28)
$ws.Cells.Item(138, 9).Value = 1951.5294  # ALC!I138 (was 1535.3096)
$ws.Cells.Item(138, 10).Value = 4121.8335  # ALC!J138 (was 3188.7068)
$ws.Cells.Item(138, 11).Value = 5854.5882  # ALC!K138 (was 4605.9288)
$ws.Cells.Item(138, 12).Value = 12365.5005  # ALC!L138 (was 9566.1204)
$ws.Cells.Item(138, 13).Value = -714.5882000000001  # ALC!M138 (was 534.0712000000003)
$ws.Cells.Item(138, 14).Value = -22645.5005  # ALC!N138 (was -19846.1204)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1378.3529  # ARM!H45 (was 1545)
$ws.Cells.Item(45, 9).Value = 1252.909  # ARM!I45 (was 1407.5)
$ws.Cells.Item(45, 10).Value = 1608.3334  # ARM!J45 (was 1682.5)
$ws.Cells.Item(45, 11).Value = 1252.909  # ARM!K45 (was 1407.5)
$ws.Cells.Item(45, 12).Value = 1608.3334  # ARM!L45 (was 1682.5)
$ws.Cells.Item(45, 13).Value = -875.9090000000001  # ARM!M45 (was -1030.5)
$ws.Cells.Item(45, 14).Value = -2362.3334  # ARM!N45 (was -2436.5)

$ws.Cells.Item(122, 8).Value = 2627.926  # ARM!H122 (was 3105.2222)
$ws.Cells.Item(122, 9).Value = 2497.4736  # ARM!I122 (was 2892.3076)
$ws.Cells.Item(122, 10).Value = 2937.75  # ARM!J122 (was 3658.8)
$ws.Cells.Item(122, 11).Value = 7492.4208  # ARM!K122 (was 8676.9228)
$ws.Cells.Item(122, 12).Value = 8813.25  # ARM!L122 (was 10976.4)
$ws.Cells.Item(122, 13).Value = -5042.4208  # ARM!M122 (was -6226.9228)
$ws.Cells.Item(122, 14).Value = -13713.25  # ARM!N122 (was -15876.4)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(63, 8).Value = 30000  # CRP!H63 (was 0)
$ws.Cells.Item(63, 10).Value = 30000  # CRP!J63 (was 0)
$ws.Cells.Item(63, 12).Value = 30000  # CRP!L63 (was 0)
$ws.Cells.Item(63, 14).Value = -31372  # CRP!N63 (was None)

$ws.Cells.Item(66, 8).Value = 30000  # CRP!H66 (was 0)
$ws.Cells.Item(66, 10).Value = 30000  # CRP!J66 (was 0)
$ws.Cells.Item(66, 12).Value = 90000  # CRP!L66 (was 0)
$ws.Cells.Item(66, 14).Value = -96864  # CRP!N66 (was None)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(106, 8).Value = 5741.2856  # CUL!H106 (was 4518.4287)
$ws.Cells.Item(106, 10).Value = 5741.2856  # CUL!J106 (was 4518.4287)
$ws.Cells.Item(106, 12).Value = 17223.8568  # CUL!L106 (was 13555.2861)
$ws.Cells.Item(106, 14).Value = -19115.8568  # CUL!N106 (was -15447.2861)

$ws.Cells.Item(131, 8).Value = 857.1134  # CUL!H131 (was 843.1)
$ws.Cells.Item(131, 9).Value = 245  # CUL!I131 (was 197.5)
$ws.Cells.Item(131, 11).Value = 735  # CUL!K131 (was 592.5)
$ws.Cells.Item(131, 13).Value = 4305  # CUL!M131 (was 4447.5)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 65215.6  # GSM!H122 (was 66911.28)
$ws.Cells.Item(122, 9).Value = 94972.59  # GSM!I122 (was 83138.19500000001)
$ws.Cells.Item(122, 10).Value = 3412.6155  # GSM!J122 (was 4032)
$ws.Cells.Item(122, 11).Value = 284917.77  # GSM!K122 (was 249414.585)
$ws.Cells.Item(122, 12).Value = 10237.8465  # GSM!L122 (was 12096)
$ws.Cells.Item(122, 13).Value = -282467.77  # GSM!M122 (was -246964.585)
$ws.Cells.Item(122, 14).Value = -15137.8465  # GSM!N122 (was -16996)

$ws.Cells.Item(123, 8).Value = 16229.467  # GSM!H123 (was 18386.223)
$ws.Cells.Item(123, 10).Value = 16229.467  # GSM!J123 (was 18386.223)
$ws.Cells.Item(123, 12).Value = 16229.467  # GSM!L123 (was 18386.223)
$ws.Cells.Item(123, 14).Value = -21129.467  # GSM!N123 (was -23286.223)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(62, 8).Value = 0  # LTW!H62 (was 29500)
$ws.Cells.Item(62, 10).Value = 0  # LTW!J62 (was 29500)
$ws.Cells.Item(62, 12).Value = 0  # LTW!L62 (was 29500)
$ws.Cells.Item(62, 14).Value = $null  # LTW!N62 (was -30748)

$ws.Cells.Item(65, 8).Value = 0  # LTW!H65 (was 29500)
$ws.Cells.Item(65, 10).Value = 0  # LTW!J65 (was 29500)
$ws.Cells.Item(65, 12).Value = 0  # LTW!L65 (was 88500)
$ws.Cells.Item(65, 14).Value = $null  # LTW!N65 (was -94740)

$ws.Cells.Item(122, 8).Value = 6946834.5  # LTW!H122 (was 6946684)
$ws.Cells.Item(122, 9).Value = 15874682  # LTW!I122 (was 10102579)
$ws.Cells.Item(122, 10).Value = 2953.3333  # LTW!J122 (was 3716)
$ws.Cells.Item(122, 11).Value = 47624046  # LTW!K122 (was 30307737)
$ws.Cells.Item(122, 12).Value = 8859.999899999999  # LTW!L122 (was 11148)
$ws.Cells.Item(122, 13).Value = -47621596  # LTW!M122 (was -30305287)
$ws.Cells.Item(122, 14).Value = -13759.9999  # LTW!N122 (was -16048)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3032.353  # WVR!H81 (was 3231.8667)
$ws.Cells.Item(81, 9).Value = 3410  # WVR!I81 (was 3994.5)
$ws.Cells.Item(81, 10).Value = 2875  # WVR!J81 (was 2954.5454)
$ws.Cells.Item(81, 11).Value = 6820  # WVR!K81 (was 7989)
$ws.Cells.Item(81, 12).Value = 5750  # WVR!L81 (was 5909.0908)
$ws.Cells.Item(81, 13).Value = -5759  # WVR!M81 (was -6928)
$ws.Cells.Item(81, 14).Value = -7872  # WVR!N81 (was -8031.0908)

$ws.Cells.Item(84, 8).Value = 3032.353  # WVR!H84 (was 3231.8667)
$ws.Cells.Item(84, 9).Value = 3410  # WVR!I84 (was 3994.5)
$ws.Cells.Item(84, 10).Value = 2875  # WVR!J84 (was 2954.5454)
$ws.Cells.Item(84, 11).Value = 34100  # WVR!K84 (was 39945)
$ws.Cells.Item(84, 12).Value = 28750  # WVR!L84 (was 29545.454)
$ws.Cells.Item(84, 13).Value = -28796  # WVR!M84 (was -34641)
$ws.Cells.Item(84, 14).Value = -39358  # WVR!N84 (was -40153.454)

$ws.Cells.Item(122, 8).Value = 57943  # WVR!H122 (was 49894)
$ws.Cells.Item(122, 9).Value = 73562.78999999999  # WVR!I122 (was 60863.47)
$ws.Cells.Item(122, 11).Value = 220688.37  # WVR!K122 (was 182590.41)
$ws.Cells.Item(122, 13).Value = -218238.37  # WVR!M122 (was -180140.41)
